$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 318, pushing existing rows 318-373 down to 320-375
$ws.Rows.Item(318).Resize(2).Insert()

# Fill new row 318
$ws.Cells.Item(318, 1).Value = 7
$ws.Cells.Item(318, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(318, 3).Value = "Ñuble"
$ws.Cells.Item(318, 4).Value = 44511
$ws.Cells.Item(318, 5).Value = 16
$ws.Cells.Item(318, 6).Value = 100112004
$ws.Cells.Item(318, 7).Value = "Cebolla"
$ws.Cells.Item(318, 8).Value = "Sin especificar"
$ws.Cells.Item(318, 9).Value = "1a nueva(o)"
$ws.Cells.Item(318, 10).Value = 10000
$ws.Cells.Item(318, 11).Value = 1100
$ws.Cells.Item(318, 12).Value = 1200
$ws.Cells.Item(318, 13).Value = 1150
$ws.Cells.Item(318, 14).Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Cells.Item(318, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(318, 16).Value = 115
$ws.Cells.Item(318, 17).Value = 10
$ws.Cells.Item(318, 18).Value = "Hortaliza"

# Fill new row 319
$ws.Cells.Item(319, 1).Value = 7
$ws.Cells.Item(319, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(319, 3).Value = "Ñuble"
$ws.Cells.Item(319, 4).Value = 44511
$ws.Cells.Item(319, 5).Value = 16
$ws.Cells.Item(319, 6).Value = 100112004
$ws.Cells.Item(319, 7).Value = "Cebolla"
$ws.Cells.Item(319, 8).Value = "Sin especificar"
$ws.Cells.Item(319, 9).Value = "1a nueva(o)"
$ws.Cells.Item(319, 10).Value = 12000
$ws.Cells.Item(319, 11).Value = 800
$ws.Cells.Item(319, 12).Value = 900
$ws.Cells.Item(319, 13).Value = 850
$ws.Cells.Item(319, 14).Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Cells.Item(319, 15).Value = "Región del Maule"
$ws.Cells.Item(319, 16).Value = 85
$ws.Cells.Item(319, 17).Value = 10
$ws.Cells.Item(319, 18).Value = "Hortaliza"
